$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.894064438964444
$ws.Range("D2").Value = 0.06674962916627414

# Row 3
$ws.Range("C3").Value = 2.366748378445679
$ws.Range("D3").Value = 0.02378143743227712
$ws.Range("G3").Value = "Sí"

# Row 4
$ws.Range("C4").Value = 1.809160665319776
$ws.Range("D4").Value = 0.07927557131157492

# Row 5
$ws.Range("C5").Value = 6.774041691908591
$ws.Range("D5").Value = [double]"8.666540662893851E-08"

# Row 6
$ws.Range("C6").Value = 0.9349617714357372
$ws.Range("D6").Value = 0.3564013547342548

# Row 7
$ws.Range("C7").Value = 0.2989310282450343
$ws.Range("D7").Value = 0.7668122523337246

# Row 8
$ws.Range("C8").Value = 4.590767148975059
$ws.Range("D8").Value = [double]"5.791053418446879E-05"

# Row 9
$ws.Range("C9").Value = -0.6728435303477377
$ws.Range("D9").Value = 0.5055945583933896

# Row 10
$ws.Range("C10").Value = 3.919927093980859
$ws.Range("D10").Value = 0.0004071148796749302

# Row 11
$ws.Range("C11").Value = 3.573263169696083
$ws.Range("D11").Value = 0.001079106529158436
